# SF User list changes - 16 May - Initial
#
# - Replace the user "Ashley Choi" with "Aadarsh Patel" on the Users sheet.
# - Make "Users" the active/selected sheet (tab), with A2 selected.
# - Update the remembered selection on the RateSheetManagement sheet to C10.

$wb = $excel.ActiveWorkbook

# Swap the user name on the Users sheet.
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Aadarsh Patel"

# Update the saved selection on RateSheetManagement.
$wsRate = $wb.Worksheets.Item("RateSheetManagement")
$wsRate.Activate()
$wsRate.Range("C10").Select()

# Leave the workbook with the Users sheet active and A2 selected.
$wsUsers.Activate()
$wsUsers.Range("A2").Select()
